$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# 1) Three tables (slides 14, 15, 16) switch from the deck's custom default
#    table style to the built-in "No Style, No Grid" table style.
# ---------------------------------------------------------------------------
$newStyleId = "{6B40F00D-F0F4-4677-A525-86D5EDB7EA5C}"

14..16 | ForEach-Object {
    $slide = $p.Slides.Item($_)
    $tableShape = $slide.Shapes.Item(1)
    if ($tableShape.HasTable) {
        $tableShape.Table.ApplyStyle($newStyleId)
    }
}

# ---------------------------------------------------------------------------
# 2) Swap the presentation's theme palette from the custom "Integral / Red
#    Violet" colors to the stock "Office" palette (dk1..folHlink, 12 slots).
# ---------------------------------------------------------------------------
$master = $p.SlideMaster
$colorScheme = $master.Theme.ThemeColorScheme

$officeColors = @(
    0,          # dk1      000000
    16777215,   # lt1      FFFFFF
    6968388,    # dk2      44546A
    15132391,   # lt2      E7E6E6
    13998939,   # accent1  5B9BD5
    3243501,    # accent2  ED7D31
    10855845,   # accent3  A5A5A5
    49407,      # accent4  FFC000
    12874308,   # accent5  4472C4
    4697456,    # accent6  70AD47
    12673797,   # hlink    0563C1
    7491477     # folHlink 954F72
)

for ($i = 1; $i -le $colorScheme.Count; $i++) {
    $colorScheme.Colors($i).RGB = $officeColors[$i - 1]
}
